# Fix the tiled texturing for the terrain.
#
# The "Draw a scene within a scene..." row (20) had its Milestone Completed
# mark (III / X) moved onto two other rows: "Substatial Use of Compute
# Shader..." (row 47) and "Apply a post process routine..." (row 75).
# Row 20 keeps its Milestone designation (III) but is no longer marked
# complete; rows 47 and 75 are now both marked complete as Milestone III.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 20: clear the "Completed" (X) mark - milestone (III) stays as-is.
$ws.Range("F20").ClearContents()

# Row 47: mark as Milestone III, Completed.
$ws.Range("E47").Value = "III"
$ws.Range("F47").Value = "X"

# Row 75: mark as Milestone III, Completed.
$ws.Range("E75").Value = "III"
$ws.Range("F75").Value = "X"

# Update the saved view/selection to match where the edit was made.
$ws.Range("F20").Select() | Out-Null
